# Update the "Forecast Comparison" sheet with correct forecast output:
# - Insert a new "Week_Start_Date" column after "Week" (shifts ASIN.. right by one)
# - Normalize the Week labels from "W01".."W16" to "W1".."W16" (no leading zero)
# - Populate the new Week_Start_Date column with each week's start date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B, shifting ASIN/MyForecast/... one column to the right
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$startDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

$dateCol = $ws.Range("B2:B17")
$dateCol.NumberFormat = "@"

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
    $ws.Cells.Item($row, 2).Value = $startDates[$i]
}
